$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from O1 onto the two
# new header cells P1 and Q1, then give them their sequential values.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Swap the I/K and M/O column values for every data row (2-25).
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# New columns P and Q are filled with 2 for every data row (2-25).
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2
